$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "No significant differences detected between overview_home_page_20240721-155359.png and overview_home_page_20240721-155448.png."
$ws.Range("B2").Value = "Success"
$ws.Range("C2").Value = "2024-07-21 15:55:10"
$ws.Range("D2").Value = "Master"
$ws.Range("E2").Value = "'0.59%"
$ws.Range("F2").Value = "overview_home_page_20240721-155359.png"
$ws.Range("G2").Value = "overview_home_page_20240721-155448.png"

# Update row 3
$ws.Range("A3").Value = "No significant differences detected between login_home_page_20240721-155402.png and login_home_page_20240721-155451.png."
$ws.Range("B3").Value = "Success"
$ws.Range("C3").Value = "2024-07-21 15:55:10"
$ws.Range("D3").Value = "Master"
$ws.Range("E3").Value = "'0.60%"
$ws.Range("F3").Value = "login_home_page_20240721-155402.png"
$ws.Range("G3").Value = "login_home_page_20240721-155451.png"

# Add new row 4
$ws.Range("A4").Value = "No significant differences detected between login_logged_in_20240721-155405.png and login_logged_in_20240721-155455.png."
$ws.Range("B4").Value = "Success"
$ws.Range("C4").Value = "2024-07-21 15:55:10"
$ws.Range("D4").Value = "Master"
$ws.Range("E4").Value = "'1.98%"
$ws.Range("F4").Value = "login_logged_in_20240721-155405.png"
$ws.Range("G4").Value = "login_logged_in_20240721-155455.png"

# Add new row 5
$ws.Range("A5").Value = "No significant differences detected between overview_displayed_20240721-155418.png and overview_displayed_20240721-155507.png."
$ws.Range("B5").Value = "Success"
$ws.Range("C5").Value = "2024-07-21 15:55:11"
$ws.Range("D5").Value = "Master"
$ws.Range("E5").Value = "'1.98%"
$ws.Range("F5").Value = "overview_displayed_20240721-155418.png"
$ws.Range("G5").Value = "overview_displayed_20240721-155507.png"
